# Equip.xlsx edit: unify the conception of DataNode, DataTable, Entity.
# The only data-level change in this revision is a rename of the sheet
# formerly called "Property1" to "DataNode" (the rest of the upstream
# diff is Excel-resave/locale noise: window geometry, revision GUIDs,
# default-font-locale swap, etc. that aren't meaningful edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet.
$ws.Name = "DataNode"

# Carry forward the author's last-active-cell selection in the frozen
# (bottomLeft) pane, same as when the workbook was resaved.
$ws.Range("L39").Select() | Out-Null
